# Updates odds values on the FlashScore weekly games sheet (rows 7, 9, 10)
# as captured by the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("G7").Value = 4.5
$ws.Range("I7").Value = 1.57
$ws.Range("J7").Value = 4.33
$ws.Range("K7").Value = 2.63
$ws.Range("M7").Value = 23
$ws.Range("N7").Value = 1.02
$ws.Range("S7").Value = 1.2
$ws.Range("T7").Value = 4.33
$ws.Range("Z7").Value = 51
$ws.Range("AB7").Value = 29
$ws.Range("AI7").Value = 11
$ws.Range("AR7").Value = 67
$ws.Range("AT7").Value = 4.33

# Row 9
$ws.Range("G9").Value = 1.57
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 5.5
$ws.Range("J9").Value = 2.1
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 2.3
$ws.Range("S9").Value = 1.29
$ws.Range("T9").Value = 3.5
$ws.Range("AC9").Value = 15
$ws.Range("AJ9").Value = 17
$ws.Range("AT9").Value = 3.5

# Row 10
$ws.Range("G10").Value = 2.1
$ws.Range("I10").Value = 3.3
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 3.75
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 11
$ws.Range("Q10").Value = 1.85
$ws.Range("R10").Value = 2
$ws.Range("W10").Value = 8.5
$ws.Range("Z10").Value = 19
$ws.Range("AC10").Value = 11
$ws.Range("AL10").Value = 26
$ws.Range("AM10").Value = 34
$ws.Range("AO10").Value = 11
$ws.Range("AW10").Value = 5.5
